$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old row 13 (the lone "3268262 - Carlos Renato Menegatti" row under
# "Docentes responsáveis:") which shifts every following row up by one.
$ws.Rows(13).Delete()

# Update the cells whose text content changed (beyond the row shift itself).
$ws.Range("B10").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("C10").Value = "3268262 - Carlos Renato Menegatti"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("C18").Value = "3268262 - Carlos Renato Menegatti"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
